# de_dg_gwq_upper.xlsx parameter update
# -------------------------------------
# Change the "base" parameterization for the
# "Division of Water Quality (SWRCB)" row (row 18) from -1 to -0.5
# across the four weighting columns C:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C18:F18").Value = -0.5

# Leave the edited range selected, matching the state the workbook was
# saved in.
$ws.Range("C18:F18").Select() | Out-Null
